# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (column E, rows 16-42) listed the 27 monthly
# periods in descending order (2003, 2002, ... 1709). This edit re-sorts
# them in ascending order (1709, 1802, ... 2003) to match the refreshed
# source database used to build the account-statement report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ascending chronological order of the 27 "Periodo Mora" values.
$periods = @(
    "1709",
    "1802", "1803", "1804", "1805", "1806", "1807", "1808", "1809", "1810", "1811", "1812",
    "1901", "1902", "1903", "1904", "1905", "1906", "1907", "1908", "1909", "1910", "1911", "1912",
    "2001", "2002", "2003"
)

$firstRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $firstRow + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
}
